# Natmi following Dr Hou advice
# Update the ligand/receptor expressing-cell counts (1 -> 3) and recompute the
# dependent expression / specificity / edge-weight columns that NATMI derives
# from them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New total expression values (sum across the now-3 expressing cells) per
# sending cluster (ligand side) and per target cluster (receptor side).
$ligandTotalBySendingCluster = @{
    "ECs"  = 31.101992
    "FAPs" = 113.744728
    "M2"   = 80.16443599999999
    "sCs"  = 40.970413
}

$receptorTotalByTargetCluster = @{
    "ECs"  = 46.076522
    "FAPs" = 151.777694
    "M2"   = 181.13147
    "sCs"  = 76.29665199999999
}

$newExpressingCells = 3

# Average expression = total expression / number of expressing cells
$ligandAvgBySendingCluster = @{}
foreach ($key in $ligandTotalBySendingCluster.Keys) {
    $ligandAvgBySendingCluster[$key] = $ligandTotalBySendingCluster[$key] / $newExpressingCells
}

$receptorAvgByTargetCluster = @{}
foreach ($key in $receptorTotalByTargetCluster.Keys) {
    $receptorAvgByTargetCluster[$key] = $receptorTotalByTargetCluster[$key] / $newExpressingCells
}

# Derived specificity = this cluster's average (or total) expression divided
# by the sum of average (or total) expression across all clusters. Since the
# expressing-cell count is now uniformly 3 for every cluster, the average-
# and total-based specificities coincide.
$ligandAvgSum = 0
foreach ($v in $ligandAvgBySendingCluster.Values) { $ligandAvgSum += $v }

$receptorAvgSum = 0
foreach ($v in $receptorAvgByTargetCluster.Values) { $receptorAvgSum += $v }

$ligandSpecBySendingCluster = @{}
foreach ($key in $ligandAvgBySendingCluster.Keys) {
    $ligandSpecBySendingCluster[$key] = $ligandAvgBySendingCluster[$key] / $ligandAvgSum
}

$receptorSpecByTargetCluster = @{}
foreach ($key in $receptorAvgByTargetCluster.Keys) {
    $receptorSpecByTargetCluster[$key] = $receptorAvgByTargetCluster[$key] / $receptorAvgSum
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp

for ($r = 2; $r -le $lastRow; $r++) {
    $sendingCluster = $ws.Cells.Item($r, 1).Value2   # column A
    $targetCluster  = $ws.Cells.Item($r, 4).Value2   # column D

    if (-not $ligandTotalBySendingCluster.ContainsKey($sendingCluster)) { continue }
    if (-not $receptorTotalByTargetCluster.ContainsKey($targetCluster)) { continue }

    $ligandTotal   = $ligandTotalBySendingCluster[$sendingCluster]
    $ligandAvg     = $ligandAvgBySendingCluster[$sendingCluster]
    $ligandSpec    = $ligandSpecBySendingCluster[$sendingCluster]

    $receptorTotal = $receptorTotalByTargetCluster[$targetCluster]
    $receptorAvg   = $receptorAvgByTargetCluster[$targetCluster]
    $receptorSpec  = $receptorSpecByTargetCluster[$targetCluster]

    $edgeAvgWeight  = $ligandAvg * $receptorAvg
    $edgeTotalWeight = $ligandTotal * $receptorTotal
    $edgeAvgSpec    = $ligandSpec * $receptorSpec
    $edgeTotalSpec  = $ligandSpec * $receptorSpec

    $ws.Cells.Item($r, 5).Value  = $newExpressingCells   # E: Ligand-expressing cells
    $ws.Cells.Item($r, 7).Value  = $ligandAvg            # G: Ligand average expression value
    $ws.Cells.Item($r, 8).Value  = $ligandTotal          # H: Ligand total expression value
    $ws.Cells.Item($r, 9).Value  = $ligandSpec           # I: Ligand derived specificity (avg)
    $ws.Cells.Item($r, 10).Value = $ligandSpec           # J: Ligand derived specificity (total)

    $ws.Cells.Item($r, 11).Value = $newExpressingCells   # K: Receptor-expressing cells
    $ws.Cells.Item($r, 13).Value = $receptorAvg          # M: Receptor average expression value
    $ws.Cells.Item($r, 14).Value = $receptorTotal        # N: Receptor total expression value
    $ws.Cells.Item($r, 15).Value = $receptorSpec         # O: Receptor derived specificity (avg)
    $ws.Cells.Item($r, 16).Value = $receptorSpec         # P: Receptor derived specificity (total)

    $ws.Cells.Item($r, 17).Value = $edgeAvgWeight        # Q: Edge average expression weight
    $ws.Cells.Item($r, 18).Value = $edgeTotalWeight       # R: Edge total expression weight
    $ws.Cells.Item($r, 19).Value = $edgeAvgSpec          # S: Edge average expression derived specificity
    $ws.Cells.Item($r, 20).Value = $edgeTotalSpec        # T: Edge total expression derived specificity
}
